$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (for the "Monk" class) - this shifts the existing
# "Rogue" column (old B) to column C, preserving its values/styles.
$ws.Columns("B").Insert()

# Fill in the Monk column (B) with its stats, row by row, mirroring the
# layout used by the Cleric (A) and Rogue (C) columns.
$ws.Range("B1").Value = "Monk"
$ws.Range("B2").Value = "1/2=2/2=3/2=4/2=5/3=6/3=7/3=8/3=9/4=10/4=11/4=12/4=13/5=14/5=15/5=16/5=17/6=18/6=19/6=20/6"
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = "None"
$ws.Range("B5").Value = "Simple Weapons=Short-sword"
$ws.Range("B6").Value = "?Artisan’s tools=?Musical Instruments"
$ws.Range("B7").Value = "Strength=Dexterity"
$ws.Range("B8").Value = "false=None"
$ws.Range("B9").Value = "1/Unarmored Defense=2/Unarmored Movement=3/Deflect Missiles=4/Slow Fall=5/Extra Attack=5/Stunning Strike=6/Ki-Empowered Strikes=7/Evasion=7/Stillness of Mind=10/Purity of Body (Immune to Disease and Poison)=13/Tongue of the Sun and Moon=14/Diamond Soul=15/Timeless Body=18/Empty Body=20/Perfect Self"
$ws.Range("B10").Value = "1/Martial Arts=2/Ki"
